$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 228, shifting rows 228:310 down to 229:311.
$ws.Rows.Item(228).Insert()

# Row 228 now has the formatting of row 229 (the old row 228, pushed down)
# copied by Insert; repopulate its content from the (now-shifted) row 229
# so the new row 228 starts as a duplicate of the old row 228, then
# overwrite the fields that actually change.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(228, $col).Value = $ws.Cells.Item(229, $col).Value()
}
$ws.Cells.Item(228, 4).Value = 44795
$ws.Cells.Item(228, 11).Value = 27000
$ws.Cells.Item(228, 12).Value = 27000
$ws.Cells.Item(228, 13).Value = 27000
$ws.Cells.Item(228, 16).Value = 450
